# dbTask.xlsx update: rename/status fixes on the "Tarea a" row, clear the
# stray blank placeholder row 4, and touch row 10 so the sheet grows to it
# (matches the A1:J9 -> A1:J10 dimension bump in the target revision).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capitalize the task title: "Tarea a" -> "Tarea A"
$ws.Range("B2").Value = "Tarea A"

# Status moved back to pending: "Finalizada" -> "En espera"
$ws.Range("D2").Value = "En espera"

# Row 4 only ever held empty placeholder cells - drop their content so the
# row collapses back down to a bare, empty row.
$ws.Range("A4:F4").ClearContents()

# Extend the sheet's used range down to row 10 (was blank/unused before)
# without leaving any cell content or formatting behind.
$ws.Rows(10).Hidden = $true
$ws.Rows(10).Hidden = $false
